$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.923.02'
$ws.Range("E2").Value = '  -1.88%  '
$ws.Range("D3").Value = '3.790.12'
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'620.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.62%  '
$ws.Range("D6").Value = "'177.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.71%  '
$ws.Range("D7").Value = '3.787.29'
$ws.Range("E7").Value = '  +2.89%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +4.26%  '
$ws.Range("E11").Value = '  -4.95%  '
$ws.Range("D12").Value = "'0.492"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("D13").Value = "'40.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("E14").Value = '  +3.17%  '
$ws.Range("D15").Value = '4.420.58'
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").Value = '3.798.80'
$ws.Range("E16").Value = '  +3.11%  '
$ws.Range("D17").Value = '69.961.65'
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("E19").Value = '  +1.29%  '
$ws.Range("D20").Value = "'16.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = "'509.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = "'9.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("D23").Value = "'0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("D24").Value = "'2.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.95%  '
$ws.Range("D25").Value = "'87.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("E26").Value = '  -1.25%  '
$ws.Range("E27").Value = '  +29.13%  '
$ws.Range("D28").Value = "'11.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").Value = "'2.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("E31").Value = '  +4.05%  '
$ws.Range("D32").Value = "'7.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.33%  '
$ws.Range("D33").Value = "'31.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  +5.93%  '
$ws.Range("D37").Value = "'6.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("E38").Value = '  +5.04%  '
$ws.Range("D39").Value = "'0.332"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.78%  '
$ws.Range("D40").Value = "'2.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").Value = "'51.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").Value = "'45.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").Value = "'419.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.81%  '
$ws.Range("D45").Value = "'2.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("D46").Value = '3.040.67'
$ws.Range("E46").Value = '  -3.95%  '
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").Value = "'27.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.71%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = "'138.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("E51").Value = '  +1.28%  '
